$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "B" and "C" rows within each year group (rows 3<->4, 7<->8, 11<->12, 15<->16)
$pairs = @(3,7,11,15)
foreach ($r1 in $pairs) {
    $r2 = $r1 + 1
    for ($col = 1; $col -le 5; $col++) {
        $v1 = $ws.Cells.Item($r1, $col).Value2
        $v2 = $ws.Cells.Item($r2, $col).Value2
        $ws.Cells.Item($r1, $col).Value2 = $v2
        $ws.Cells.Item($r2, $col).Value2 = $v1
    }
}

# Delete columns F and G (产销率 and 销售量 instantaneous columns)
$ws.Range("F1:G17").Delete()
